# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.75"
$ws.Range("E2").Value = "'0.62%"
$ws.Range("E3").Value = "'1.02%"
$ws.Range("D4").Value = "'4.834"
$ws.Range("E4").Value = "'-0.79%"
$ws.Range("D5").Value = "'0.06390"
$ws.Range("E5").Value = "'0.33%"
$ws.Range("D6").Value = "'7.045"
$ws.Range("E6").Value = "'0.97%"
$ws.Range("D7").Value = "'1.312"
$ws.Range("E7").Value = "'4.89%"
$ws.Range("D8").Value = "'0.8953"
$ws.Range("E8").Value = "'1.56%"
$ws.Range("D9").Value = "'0.1539"
$ws.Range("E9").Value = "'1.71%"
$ws.Range("D10").Value = "'0.06887"
$ws.Range("E10").Value = "'34.36%"
$ws.Range("D11").Value = "'0.07537"
$ws.Range("E11").Value = "'0.11%"
$ws.Range("D12").Value = "'0.02950"
$ws.Range("E12").Value = "'-0.59%"
$ws.Range("D13").Value = "'0.09002"
$ws.Range("E13").Value = "'-0.23%"
$ws.Range("D14").Value = "'0.001565"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("D15").Value = "'0.0006512"
$ws.Range("E15").Value = "'1.14%"
$ws.Range("D16").Value = "'0.006012"
$ws.Range("E16").Value = "'1.73%"
$ws.Range("D17").Value = "'3.486"
$ws.Range("E17").Value = "'0.61%"
$ws.Range("D18").Value = "'3.325"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("E19").Value = "'-1.91%"
$ws.Range("E20").Value = "'-1.26%"
$ws.Range("E21").Value = "'1.02%"
$ws.Range("D22").Value = "'3.897"
$ws.Range("E22").Value = "'-0.34%"
$ws.Range("D23").Value = "'0.04423"
$ws.Range("E23").Value = "'0.06%"
$ws.Range("E25").Value = "'0.27%"
$ws.Range("D26").Value = "'0.004277"
$ws.Range("E26").Value = "'10.22%"
$ws.Range("D28").Value = "'0.0001179"
$ws.Range("E28").Value = "'-1.72%"
$ws.Range("E29").Value = "'-14.58%"
$ws.Range("D40").Value = "'0.04063"
$ws.Range("E40").Value = "'-2.41%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1412"
$ws.Range("E41").Value = "'19.53%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006591"
$ws.Range("E42").Value = "'-4.06%"
$ws.Range("D43").Value = "'0.002078"
$ws.Range("E43").Value = "'-0.53%"
$ws.Range("D44").Value = "'0.01104"
$ws.Range("E44").Value = "'-1.63%"
$ws.Range("D45").Value = "'0.00005555"
$ws.Range("E45").Value = "'7.16%"
$ws.Range("D46").Value = "'1.561"
$ws.Range("E46").Value = "'5.01%"
$ws.Range("E47").Value = "'-8.69%"
